# Applies the 25-May-2023 GitHub Actions "Updated cryptos list" refresh:
# - D/E columns (Price, Volume(1h)) get refreshed quotes for rows 2-51.
# - Rows 34/35 (Filecoin / RenderToken) swap places with new figures, since
#   the coin ranking API reordered them that day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these cells hold plain text (e.g. "26.533.94", "  -0.17%  ") rather
# than real numbers -- the source sheet uses dots as thousands separators and
# keeps the percent sign + padding baked into the string. Assigning a bare
# numeric-looking string via .Value lets Excel "helpfully" reinterpret it as a
# number/percentage, so we briefly force Text format, assign, then drop back to
# the default "Normal" style (matching every other cell in these columns, which
# carry no explicit style) so only the cell *contents* change.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$updates = @(
    @{ Cell = 'D2'; Text = '26.533.94' }
    @{ Cell = 'E2'; Text = '  -0.17%  ' }
    @{ Cell = 'D3'; Text = '1.813.42' }
    @{ Cell = 'E3'; Text = '  -0.24%  ' }
    @{ Cell = 'E4'; Text = '  -0.38%  ' }
    @{ Cell = 'E5'; Text = '  -0.37%  ' }
    @{ Cell = 'E6'; Text = '  -1.02%  ' }
    @{ Cell = 'E7'; Text = '  -0.85%  ' }
    @{ Cell = 'D8'; Text = '0.3588' }
    @{ Cell = 'E8'; Text = '  -2.23%  ' }
    @{ Cell = 'D9'; Text = '46.39' }
    @{ Cell = 'E9'; Text = '  +2.63%  ' }
    @{ Cell = 'D10'; Text = '0.07104' }
    @{ Cell = 'E10'; Text = '  -0.69%  ' }
    @{ Cell = 'D11'; Text = '0.8901' }
    @{ Cell = 'E11'; Text = '  +1.33%  ' }
    @{ Cell = 'D12'; Text = '0.07740' }
    @{ Cell = 'E12'; Text = '  -0.55%  ' }
    @{ Cell = 'D13'; Text = '19.29' }
    @{ Cell = 'E13'; Text = '  -0.48%  ' }
    @{ Cell = 'D14'; Text = '1.823.33' }
    @{ Cell = 'E14'; Text = '  +0.50%  ' }
    @{ Cell = 'E15'; Text = '  -0.70%  ' }
    @{ Cell = 'D16'; Text = '6.294' }
    @{ Cell = 'E16'; Text = '  -1.33%  ' }
    @{ Cell = 'D17'; Text = '85.68' }
    @{ Cell = 'E17'; Text = '  -0.83%  ' }
    @{ Cell = 'D18'; Text = '1.006' }
    @{ Cell = 'E18'; Text = '  -0.45%  ' }
    @{ Cell = 'D19'; Text = '0.000008495' }
    @{ Cell = 'E19'; Text = '  -1.50%  ' }
    @{ Cell = 'E20'; Text = '  -0.41%  ' }
    @{ Cell = 'D21'; Text = '26.576.88' }
    @{ Cell = 'E21'; Text = '  -0.28%  ' }
    @{ Cell = 'D22'; Text = '14.11' }
    @{ Cell = 'E22'; Text = '  -1.18%  ' }
    @{ Cell = 'D23'; Text = '4.950' }
    @{ Cell = 'E23'; Text = '  -1.19%  ' }
    @{ Cell = 'E24'; Text = '  +0.17%  ' }
    @{ Cell = 'D25'; Text = '1.940' }
    @{ Cell = 'E25'; Text = '  -2.51%  ' }
    @{ Cell = 'D26'; Text = '151.91' }
    @{ Cell = 'E26'; Text = '  +0.18%  ' }
    @{ Cell = 'D27'; Text = '17.79' }
    @{ Cell = 'E27'; Text = '  -1.04%  ' }
    @{ Cell = 'D28'; Text = '2.020' }
    @{ Cell = 'E28'; Text = '  -3.06%  ' }
    @{ Cell = 'D29'; Text = '112.39' }
    @{ Cell = 'E29'; Text = '  -0.71%  ' }
    @{ Cell = 'D30'; Text = '4.827' }
    @{ Cell = 'E30'; Text = '  -0.86%  ' }
    @{ Cell = 'D31'; Text = '0.08708' }
    @{ Cell = 'E31'; Text = '  +0.13%  ' }
    @{ Cell = 'D32'; Text = '3.137' }
    @{ Cell = 'E32'; Text = '  +2.22%  ' }
    @{ Cell = 'D33'; Text = '0.7404' }
    @{ Cell = 'E33'; Text = '  +0.39%  ' }
    @{ Cell = 'B34'; Text = 'Filecoin' }
    @{ Cell = 'C34'; Text = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D34'; Text = '4.419' }
    @{ Cell = 'E34'; Text = '  -2.48%  ' }
    @{ Cell = 'B35'; Text = 'RenderToken' }
    @{ Cell = 'C35'; Text = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D35'; Text = '2.710' }
    @{ Cell = 'E35'; Text = '  -0.27%  ' }
    @{ Cell = 'E36'; Text = '  -1.09%  ' }
    @{ Cell = 'E37'; Text = '  -1.34%  ' }
    @{ Cell = 'E38'; Text = '  -0.62%  ' }
    @{ Cell = 'D39'; Text = '2.915' }
    @{ Cell = 'E39'; Text = '  +0.21%  ' }
    @{ Cell = 'D40'; Text = '0.05080' }
    @{ Cell = 'E40'; Text = '  -0.81%  ' }
    @{ Cell = 'D41'; Text = '0.5102' }
    @{ Cell = 'E41'; Text = '  +1.35%  ' }
    @{ Cell = 'D42'; Text = '6.769' }
    @{ Cell = 'E42'; Text = '  -3.65%  ' }
    @{ Cell = 'E43'; Text = '  -3.40%  ' }
    @{ Cell = 'D44'; Text = '8.028' }
    @{ Cell = 'E44'; Text = '  -2.18%  ' }
    @{ Cell = 'D45'; Text = '0.4698' }
    @{ Cell = 'E45'; Text = '  +1.37%  ' }
    @{ Cell = 'E46'; Text = '  -0.43%  ' }
    @{ Cell = 'D47'; Text = '9.981' }
    @{ Cell = 'E47'; Text = '  -0.60%  ' }
    @{ Cell = 'D48'; Text = '98.73' }
    @{ Cell = 'E48'; Text = '  -2.48%  ' }
    @{ Cell = 'E49'; Text = '  -1.97%  ' }
    @{ Cell = 'E50'; Text = '  -0.47%  ' }
    @{ Cell = 'D51'; Text = '63.78' }
    @{ Cell = 'E51'; Text = '  -1.32%  ' }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Range($u.Cell) $u.Text
}
